$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Ziekenhuisorganisaties", "Status inleveren routekaart"),
    @("Albert Schweitzer ziekenhuis", "voorlopig"),
    @("Alrijne Zorggroep", "voorlopig"),
    @("Amphia Ziekenhuis", "voorlopig"),
    @("Amsterdam UMC", "voorlopig"),
    @("Antonius Zorggroep", "definitief en/of vastgesteld RvB"),
    @("Bravis ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Canisius-Wilhelmina Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Catharina Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Deventer Ziekenhuis", "voorlopig"),
    @("Elisabeth-TweeSteden Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Elkerliek Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Erasmus MC", "voorlopig"),
    @("Flevoziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Franciscus Gasthuis & Vlietland Groep", "definitief en/of vastgesteld RvB"),
    @("Gelre Ziekenhuizen", "definitief en/of vastgesteld RvB"),
    @("Groene Hart Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("IJsselland Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Isala Klinieken", "definitief en/of vastgesteld RvB"),
    @("Jeroen Bosch Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("LUMC", "voorlopig"),
    @("Laurentius Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Maasstad ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Maastricht UMC+", "definitief en/of vastgesteld RvB"),
    @("Martini Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Medisch Centrum Leeuwarden", "definitief en/of vastgesteld RvB"),
    @("Medisch Spectrum Twente", "voorlopig"),
    @("Máxima Medisch Centrum", "definitief en/of vastgesteld RvB"),
    @("Nij Smellinghe", "voorlopig"),
    @("Noordwest Ziekenhuisgroep", "definitief en/of vastgesteld RvB"),
    @("Ommelander Ziekenhuis Groep", "definitief en/of vastgesteld RvB"),
    @("Onze Lieve Vrouwe Gasthuis", "definitief en/of vastgesteld RvB"),
    @("Pantein", "definitief en/of vastgesteld RvB"),
    @("Prinses Máxima Centrum", "definitief en/of vastgesteld RvB"),
    @("Radboudumc", "voorlopig"),
    @("Rode Kruis Ziekenhuis", "voorlopig"),
    @("Slingeland Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Spaarne Gasthuis", "definitief en/of vastgesteld RvB"),
    @("St Jansdal", "definitief en/of vastgesteld RvB"),
    @("St. Antonius Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("St. Jans Gasthuis", "definitief en/of vastgesteld RvB"),
    @("Stichting BovenIJ ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Stichting Dijklander Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Stichting Reinier Haga Groep", "definitief en/of vastgesteld RvB"),
    @("Stichting Reinier de Graaf Groep", "definitief en/of vastgesteld RvB"),
    @("Stichting Rijnstate Ziekenhuis", "definitief en/of vastgesteld RvB"),
    @("Stichting St. Anna Zorggroep", "definitief en/of vastgesteld RvB"),
    @("Stichting Tergooi", "definitief en/of vastgesteld RvB"),
    @("Stichting VieCuri-Vitaal", "definitief en/of vastgesteld RvB"),
    @("Stichting Ziekenhuisgroep Twente (ZGT)", "voorlopig"),
    @("Tjongerschans", "definitief en/of vastgesteld RvB"),
    @("Treant Zorggroep", "voorlopig"),
    @("UMC Utrecht", "definitief en/of vastgesteld RvB"),
    @("UMCG", "voorlopig"),
    @("Wilhelmina Ziekenhuis Assen", "definitief en/of vastgesteld RvB"),
    @("Ziekenhuis Amstelland", "voorlopig"),
    @("Ziekenhuis Bernhoven", "definitief en/of vastgesteld RvB"),
    @("Ziekenhuis De Gelderse Vallei", "definitief en/of vastgesteld RvB"),
    @("Zuyderland Ziekenhuis", "definitief en/of vastgesteld RvB")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $ws.Cells.Item($i + 1, 1).Value = $row[0]
    $ws.Cells.Item($i + 1, 2).Value = $row[1]
}
